# Adds 3 new upcoming-fixture rows (116-118) to the Canada Premier League
# sheet, mirroring the source feed update described in the commit message
# ("Atualizacao de bases das ligas" - league base update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append right after the current last row (115).
$rows = @(
    @{
        Row = 116; A = 114; B = "7802952"; C = "Canada Premier League";
        D = 45464.91666666666; E = "Cavalry FC"; F = "Atletico Ottawa";
        L = 2.15; M = 3.1; N = 3.2; O = 1.95; P = 3.2; Q = 3.6; R = -0.5;
        S = 2.025; T = 1.775; U = 2.25; V = 1.9; W = 1.9; X = 0; Y = 0; Z = 0
    },
    @{
        Row = 117; A = 115; B = "7802884"; C = "Canada Premier League";
        D = 45466.70833333334; E = "Forge FC"; F = "Valour FC";
        L = 1.444; M = 4; N = 6.5; O = 1.533; P = 4; Q = 5.5; R = -1;
        S = 1.825; T = 1.975; U = 2.75; V = 1.975; W = 1.825; X = 0; Y = 0; Z = 0
    },
    @{
        Row = 118; A = 116; B = "7802953"; C = "Canada Premier League";
        D = 45466.83333333334; E = "Vancouver FC"; F = "HFX Wanderers";
        L = 2.3; M = 3.4; N = 2.7; O = 2.25; P = 3.4; Q = 2.7; R = -0.25;
        S = 2.025; T = 1.775; U = 2.5; V = 1.85; W = 1.95; X = 0; Y = 0; Z = 0
    }
)

foreach ($r in $rows) {
    $row = $r.Row

    # Column A uses the bordered/bold "index" style already applied to A2:A115.
    $ws.Range("A115").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)
    $ws.Range("A$row").Value = $r.A

    # Column B ("id") holds large numeric-looking ids that are stored as TEXT
    # in the source feed (same as the new shared strings "7802952"/"7802884"/
    # "7802953"). Force text storage, write the value, then strip the
    # temporary "@" number-format back to the sheet's normal/general style.
    $ws.Range("B$row").NumberFormat = "@"
    $ws.Range("B$row").Value = $r.B

    $ws.Range("C$row").Value = $r.C

    # Column D uses the custom YYYY-MM-DD HH:MM:SS date/time style.
    $ws.Range("D115").Copy()
    $ws.Range("D$row").PasteSpecial(-4122)
    $ws.Range("D$row").Value = $r.D

    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F

    # G/H/I/J/K (scores + result) are left blank: these are upcoming fixtures.

    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = $r.O
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("R$row").Value = $r.R
    $ws.Range("S$row").Value = $r.S
    $ws.Range("T$row").Value = $r.T
    $ws.Range("U$row").Value = $r.U
    $ws.Range("V$row").Value = $r.V
    $ws.Range("W$row").Value = $r.W
    $ws.Range("X$row").Value = $r.X
    $ws.Range("Y$row").Value = $r.Y
    $ws.Range("Z$row").Value = $r.Z
}

# Strip the transient "@" text format from column B so the cells keep the
# sheet's default (unstyled) look, matching every other "id" column that
# happens to be plain numeric - only the string TYPE needs to stick, not a
# visible text format.
$ws.Range("B2").Copy()
$ws.Range("B116:B118").PasteSpecial(-4122)
